# ---------------------------------------------------------------------------
# Commit: "set index to Room ID in final df"
#
# 1) "readme" sheet: the metadata table's columns were reordered from
#      index | Author | sheet_name | Date | JobNo
#    to
#      index | Date | JobNo | sheet_name | Author
#    (the table header + every data row moves together).
#
# 2) "Project Information" sheet: the "Date of Analysis" value cell was
#    refreshed to a later timestamp from the same run.
#
# 3) Each "Results, Air Speed *" sheet (and its ListObject/table): the first
#    two columns were swapped so "Room ID" becomes the index/first column and
#    "Room Name" becomes the second column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) readme sheet - reorder Author / sheet_name / Date / JobNo -> Date / JobNo / sheet_name / Author
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("readme")

$sheetNames = @(
    "Project Information",
    "Criterion Definitions",
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

# Header row
$readme.Range("B1").Value = "Date"
$readme.Range("C1").Value = "JobNo"
$readme.Range("D1").Value = "sheet_name"
$readme.Range("E1").Value = "Author"

# Data rows: Date and JobNo and Author are constant for every row;
# sheet_name varies per row (one row per worksheet in the workbook).
for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $row = $i + 2
    $readme.Range("B$row").Value = "20220308"
    $readme.Range("C$row").Value = "/c/e"
    $readme.Range("D$row").Value = $sheetNames[$i]
    $readme.Range("E$row").Value = "jovyan"
}

# ---------------------------------------------------------------------------
# 2) Project Information sheet - refresh the analysis timestamp
# ---------------------------------------------------------------------------
$projInfo = $wb.Worksheets.Item("Project Information")
$projInfo.Range("B11").Value = "2022-03-08 14:52:45.801563"

# ---------------------------------------------------------------------------
# 3) Results sheets - swap "Room Name" (col A) and "Room ID" (col B)
# ---------------------------------------------------------------------------
$roomNames = @(
    "A_01_XX_XX_ApartmentSW",
    "A_01_XX_XX_ApartmentS1",
    "A_01_XX_XX_ApartmentS2",
    "A_01_XX_XX_ApartmentN1",
    "A_01_XX_XX_ApartmentNW",
    "A_01_XX_XX_ApartmentNE",
    "A_01_XX_XX_ApartmentN2",
    "A_02_XX_XX_ApartmentSW",
    "A_02_XX_XX_ApartmentS1",
    "A_02_XX_XX_ApartmentSE",
    "A_02_XX_XX_ApartmentS2",
    "A_02_XX_XX_ApartmentN1",
    "A_02_XX_XX_ApartmentNW",
    "A_02_XX_XX_ApartmentNE",
    "A_02_XX_XX_ApartmentN2",
    "A_03_XX_XX_ApartmentSW",
    "A_03_XX_XX_ApartmentS1",
    "A_03_XX_XX_ApartmentSE",
    "A_03_XX_XX_ApartmentS2",
    "A_03_XX_XX_ApartmentN1",
    "A_03_XX_XX_ApartmentNW",
    "A_03_XX_XX_ApartmentNE",
    "A_03_XX_XX_ApartmentN2",
    "A_04_XX_XX_ApartmentSW",
    "A_04_XX_XX_ApartmentS1",
    "A_04_XX_XX_ApartmentSE",
    "A_04_XX_XX_ApartmentS2",
    "A_04_XX_XX_ApartmentN1",
    "A_04_XX_XX_ApartmentNW",
    "A_04_XX_XX_ApartmentNE",
    "A_04_XX_XX_ApartmentN2"
)

$roomIds = @(
    "1S000001","1S000002","1S000003","1S000004","1S000005","1S000006","1S000007",
    "2N000001","2N000002","2N000003","2N000004","2N000005","2N000006","2N000007","2N000008",
    "3R000001","3R000002","3R000003","3R000004","3R000005","3R000006","3R000007","3R000008",
    "4T000001","4T000002","4T000003","4T000004","4T000005","4T000006","4T000007","4T000008"
)

$resultSheets = @(
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

foreach ($sheetName in $resultSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header: A=Room ID, B=Room Name
    $ws.Range("A1").Value = "Room ID"
    $ws.Range("B1").Value = "Room Name"

    for ($i = 0; $i -lt $roomNames.Count; $i++) {
        $row = $i + 2
        $ws.Range("A$row").Value = $roomIds[$i]
        $ws.Range("B$row").Value = $roomNames[$i]
    }
}
